$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I15").Value = 0.535526637811788
$ws.Range("H16").Value = 0.5766911554241068
$ws.Range("G17").Value = 0.6272238950261231
$ws.Range("F18").Value = 0.6666911554241067
$ws.Range("E19").Value = 0.6966911554241066
$ws.Range("D20").Value = 0.4271648845785767
$ws.Range("C21").Value = 0.4775315349050862
$ws.Range("B22").Value = 0.32386998960715
